$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 (quarter 01-01-2021) with revised figures ---
$ws.Range("F74").Value  = -51
$ws.Range("P74").Value  = 484
$ws.Range("R74").Value  = 1243
$ws.Range("T74").Value  = -5
$ws.Range("U74").Value  = 2174
$ws.Range("V74").Value  = 666
$ws.Range("W74").Value  = 1736
$ws.Range("Y74").Value  = -214
$ws.Range("Z74").Value  = 10441
$ws.Range("AA74").Value = 7732
$ws.Range("AB74").Value = 1326
$ws.Range("AD74").Value = 1633

# --- Append new row 75 (quarter 01-04-2021) ---
# Column A holds a date-like label ("01-04-2021") that must be stored as
# plain text (shared string), matching the rest of column A. Temporarily
# force a text number format so Excel's automatic date recognition does not
# convert the string into a date serial value, then restore the default
# style so the cell keeps no explicit style (same as its neighbours).
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").Style = "Normal"

$ws.Range("B75").Value  = -3990
$ws.Range("C75").Value  = -4109
$ws.Range("D75").Value  = 0
$ws.Range("E75").Value  = 118
$ws.Range("F75").Value  = -169
$ws.Range("G75").Value  = 4523
$ws.Range("H75").Value  = 4446
$ws.Range("I75").Value  = 0
$ws.Range("J75").Value  = 77
$ws.Range("K75").Value  = 701
$ws.Range("L75").Value  = 5
$ws.Range("M75").Value  = 696
$ws.Range("N75").Value  = -1909
$ws.Range("O75").Value  = 1909
$ws.Range("P75").Value  = -7115
$ws.Range("Q75").Value  = 0
$ws.Range("R75").Value  = -6609
$ws.Range("S75").Value  = -614
$ws.Range("T75").Value  = 108
$ws.Range("U75").Value  = 1722
$ws.Range("V75").Value  = 241
$ws.Range("W75").Value  = 1495
$ws.Range("X75").Value  = -33
$ws.Range("Y75").Value  = 18
$ws.Range("Z75").Value  = 6651
$ws.Range("AA75").Value = 1203
$ws.Range("AB75").Value = 3170
$ws.Range("AC75").Value = -176
$ws.Range("AD75").Value = 2453
